$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.472.64"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "'3.081.14"
$ws.Range("E3").Value = "  -0.43%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'545.39"
$ws.Range("E5").Value = "  -0.47%  "

$ws.Range("D6").Value = "'139.33"
$ws.Range("E6").Value = "  +1.81%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").Value = "'3.073.98"
$ws.Range("E8").Value = "  -0.36%  "

$ws.Range("D9").Value = "'0.501"
$ws.Range("E9").Value = "  +0.85%  "

$ws.Range("E10").Value = "  -0.60%  "

$ws.Range("D11").Value = "'6.37"
$ws.Range("E11").Value = "  +1.55%  "

$ws.Range("D12").Value = "'0.457"
$ws.Range("E12").Value = "  -2.96%  "

$ws.Range("D13").Value = "'35.03"
$ws.Range("E13").Value = "  -1.35%  "

$ws.Range("D14").Value = "'0.0000224"
$ws.Range("E14").Value = "  +3.35%  "

$ws.Range("D15").Value = "'3.583.89"
$ws.Range("E15").Value = "  -0.30%  "

$ws.Range("D16").Value = "'63.484.47"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("D18").Value = "'3.081.78"
$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("E19").Value = "  -1.37%  "

$ws.Range("D20").Value = "'475.39"
$ws.Range("E20").Value = "  -2.90%  "

$ws.Range("D21").Value = "'13.49"
$ws.Range("E21").Value = "  -1.13%  "

$ws.Range("D22").Value = "'0.701"
$ws.Range("E22").Value = "  -2.44%  "

$ws.Range("D23").Value = "'7.09"
$ws.Range("E23").Value = "  -1.99%  "

$ws.Range("D24").Value = "'78.74"
$ws.Range("E24").Value = "  -0.51%  "

$ws.Range("D25").Value = "'12.23"
$ws.Range("E25").Value = "  -1.26%  "

$ws.Range("E26").Value = "  +0.20%  "

$ws.Range("D27").Value = "'2.73"
$ws.Range("E27").Value = "  -1.04%  "

$ws.Range("D28").Value = "'7.99"
$ws.Range("E28").Value = "  -5.66%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").Value = "'26.22"
$ws.Range("E30").Value = "  -1.44%  "

$ws.Range("E31").Value = "  -3.53%  "

$ws.Range("E32").Value = "  +2.64%  "

$ws.Range("D33").Value = "'58.00"
$ws.Range("E33").Value = "  -0.16%  "

$ws.Range("D34").Value = "'2.32"
$ws.Range("E34").Value = "  -7.45%  "

$ws.Range("D35").Value = "'5.45"
$ws.Range("E35").Value = "  +6.23%  "

$ws.Range("D36").Value = "'494.11"
$ws.Range("E36").Value = "  -3.89%  "

$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("D38").Value = "'3.263.28"
$ws.Range("E38").Value = "  +3.38%  "

$ws.Range("E39").Value = "  +0.26%  "

$ws.Range("D40").Value = "'0.0799"
$ws.Range("E40").Value = "  -0.43%  "

$ws.Range("E41").Value = "  -1.53%  "

$ws.Range("D42").Value = "'8.15"
$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").Value = "'2.62"
$ws.Range("E43").Value = "  -1.57%  "

$ws.Range("E44").Value = "  -2.00%  "

$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'25.44"
$ws.Range("E46").Value = "  +0.86%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'123.00"
$ws.Range("E47").Value = "  +1.62%  "

$ws.Range("E48").Value = "  -0.93%  "

$ws.Range("D49").Value = "'0.0₃0531"
$ws.Range("E49").Value = "  +5.49%  "

$ws.Range("E50").Value = "  +0.79%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.03"
$ws.Range("E51").Value = "  -0.30%  "
